$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.048.46"
$ws.Range("E2").Value = "  +0.34%  "

# Row 3
$ws.Range("D3").Value = "3.581.86"
$ws.Range("E3").Value = "  +2.15%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "604.31"
$ws.Range("E5").Value = "  -0.42%  "

# Row 6
$ws.Range("D6").Value = "195.73"
$ws.Range("E6").Value = "  -1.60%  "

# Row 7
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  -0.20%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").Value = "0.205"
$ws.Range("E9").Value = "  -3.02%  "

# Row 10
$ws.Range("E10").Value = "  -0.84%  "

# Row 11
$ws.Range("D11").Value = "53.83"
$ws.Range("E11").Value = "  -0.80%  "

# Row 12
$ws.Range("E12").Value = "  -0.28%  "

# Row 13
$ws.Range("D13").Value = "9.54"
$ws.Range("E13").Value = "  -1.01%  "

# Row 14
$ws.Range("D14").Value = "4.139.87"
$ws.Range("E14").Value = "  +1.76%  "

# Row 15
$ws.Range("D15").Value = "598.34"
$ws.Range("E15").Value = "  +0.13%  "

# Row 16
$ws.Range("E16").Value = "  +2.37%  "

# Row 17
$ws.Range("D17").Value = "19.27"
$ws.Range("E17").Value = "  +1.35%  "

# Row 18
$ws.Range("D18").Value = "70.223.55"
$ws.Range("E18").Value = "  +0.35%  "

# Row 19
$ws.Range("D19").Value = "3.558.50"
$ws.Range("E19").Value = "  +1.13%  "

# Row 20
$ws.Range("E20").Value = "  +1.78%  "

# Row 21
$ws.Range("D21").Value = "0.993"
$ws.Range("E21").Value = "  -0.51%  "

# Row 22
$ws.Range("D22").Value = "17.85"
$ws.Range("E22").Value = "  +0.45%  "

# Row 23
$ws.Range("D23").Value = "5.14"
$ws.Range("E23").Value = "  +1.53%  "

# Row 24
$ws.Range("D24").Value = "102.64"
$ws.Range("E24").Value = "  -1.40%  "

# Row 25
$ws.Range("E25").Value = "  -0.26%  "

# Row 26
$ws.Range("E26").Value = "  -1.67%  "

# Row 27
$ws.Range("D27").Value = "10.81"
$ws.Range("E27").Value = "  -1.66%  "

# Row 28
$ws.Range("D28").Value = "9.58"
$ws.Range("E28").Value = "  -2.41%  "

# Row 29
$ws.Range("D29").Value = "33.82"
$ws.Range("E29").Value = "  -0.15%  "

# Row 30
$ws.Range("D30").Value = "7.13"
$ws.Range("E30").Value = "  -1.31%  "

# Row 31
$ws.Range("D31").Value = "4.28"
$ws.Range("E31").Value = "  -6.42%  "

# Row 32
$ws.Range("D32").Value = "12.34"
$ws.Range("E32").Value = "  -3.37%  "

# Row 33
$ws.Range("E33").Value = "  -0.11%  "

# Row 34
$ws.Range("D34").Value = "63.39"
$ws.Range("E34").Value = "  -0.70%  "

# Row 35
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "3.860.19"
$ws.Range("E35").Value = "  +3.71%  "

# Row 36
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "3.23"
$ws.Range("E36").Value = "  +6.81%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0822"
$ws.Range("E37").Value = "  +3.33%  "

# Row 38
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "523.91"
$ws.Range("E38").Value = "  +0.39%  "

# Row 39
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.16%  "

# Row 40
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "37.12"
$ws.Range("E40").Value = "  +0.64%  "

# Row 41
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "0.393"
$ws.Range("E41").Value = "  +0.24%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "3.60"
$ws.Range("E42").Value = "  +1.28%  "

# Row 43
$ws.Range("E43").Value = "  -2.00%  "

# Row 44
$ws.Range("E44").Value = "  -1.92%  "

# Row 45
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "2.86"
$ws.Range("E45").Value = "  -0.32%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.140"
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
$ws.Range("D47").Value = "3.29"
$ws.Range("E47").Value = "  -0.45%  "

# Row 48
$ws.Range("E48").Value = "  -2.25%  "

# Row 49
$ws.Range("E49").Value = "  -0.02%  "

# Row 50
$ws.Range("D50").Value = "0.000250"
$ws.Range("E50").Value = "  +4.30%  "

# Row 51
$ws.Range("D51").Value = "1.30"
$ws.Range("E51").Value = "  +1.99%  "
